$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark currently on the title paragraph.
$d.Bookmarks("_GoBack").Delete()

# 2. Append new content at the very end of the document body.
#    Paragraph: "Role Assigning: Add" / "/Remove users from roles" (sz 36)
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Font.Size = 18
$r.Font.SizeBi = 18
$r.InsertAfter("Role Assigning: Add")
$r.Collapse(0)
$r.Font.Size = 18
$r.Font.SizeBi = 18
$r.InsertAfter("/Remove users from roles")

# Paragraph: "Admins are able to add or remove a user from a role. " (sz 28)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Font.Size = 14
$r.Font.SizeBi = 14
$r.InsertAfter("Admins are able to add or remove a user from a role. ")

# Paragraph: "Projects" (sz 36)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Font.Size = 18
$r.Font.SizeBi = 18
$r.InsertAfter("Projects")

# Paragraph: "PMs and Admins can create Projects then can assign developers to the project." (sz 28)
# followed by the _GoBack bookmark.
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Font.Size = 14
$r.Font.SizeBi = 14
$r.InsertAfter("PMs and Admins can create Projects then can assign developers to the project.")

# Paragraph: empty trailing paragraph (sz 28)
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.Font.Size = 14
$r.Font.SizeBi = 14

# Now place the _GoBack bookmark at the end of the "PMs and Admins..." paragraph
# (immediately after its text, before the paragraph mark).
$bmPara = $d.Paragraphs.Last.Previous
$bmRange = $bmPara.Range
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "edit complete"
